$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell K45 from FALSE to TRUE
$ws.Cells.Item(45, 11).Value = $true

# Row data to append (rows 53-55) for 2025-02-18
$rows = @(
    @("2025-02-18", "sleep", $true, $false, $true, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true),
    @("2025-02-18", "activity", $false, $false, $true, $false, $true, $true, $false, $false, $false, $false, $true, $false, $false),
    @("2025-02-18", "weekly_activity", $false, $false, $true, $false, $false, $false, $false, $true, $false, $false, $false, $false, $false)
)

$startRow = 53
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    # Force column A to be stored as text (matches existing rows which store the date as a string),
    # then reset the cell style so no stray number-format style sticks around.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
    $ws.Cells.Item($r, 1).Style = "Normal"
}
